# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps for the
# ae9ac615-... handback row (row 4) on both locale sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-12 00:30:25"
$wsZhCn.Range("H4").Value = "2016-03-12 00:30:42"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-12 00:30:28"
$wsDeDe.Range("H4").Value = "2016-03-12 00:30:47"
